# Apply updates to the "by Coach" sheet:
#  - Toggle "Yes"/"No" values in column C for a set of rows
#  - Update the frozen-pane top-left cell and the current selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

# Rows whose column C "Started" value flips between Yes and No
$rows = 3,4,5,6,7,10,18,20,22,23,32,37,50,55,64,65,67,69,78,82

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "Yes") {
        $cell.Value2 = "No"
    } elseif ($cell.Value2 -eq "No") {
        $cell.Value2 = "Yes"
    }
}

# Update the current selection to reflect the scrolled view (row 1 stays
# frozen; only the visible/selected cell moves down to C83)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 65
$ws.Range("C83").Select()
